$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily mark the Price/Volume columns as Text so that numeric-looking
# strings (e.g. "0.999", "324.15") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '47.370.48'
$ws.Range('E2').Value = '  +3.20%  '
$ws.Range('D3').Value = '2.508.97'
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '324.15'
$ws.Range('E5').Value = '  +0.89%  '
$ws.Range('D6').Value = '110.06'
$ws.Range('E6').Value = '  +5.59%  '
$ws.Range('E7').Value = '  +1.39%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.538'
$ws.Range('E9').Value = '  +1.15%  '
$ws.Range('D10').Value = '39.34'
$ws.Range('E10').Value = '  +9.94%  '
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('D13').Value = '18.61'
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('E14').Value = '  +2.72%  '
$ws.Range('D15').Value = '2.901.06'
$ws.Range('E15').Value = '  +2.52%  '
$ws.Range('D16').Value = '2.512.40'
$ws.Range('E16').Value = '  +2.72%  '
$ws.Range('E17').Value = '  +2.72%  '
$ws.Range('D18').Value = '47.319.97'
$ws.Range('E18').Value = '  +3.38%  '
$ws.Range('D19').Value = '12.93'
$ws.Range('E19').Value = '  +3.30%  '
$ws.Range('E20').Value = '  +4.66%  '
$ws.Range('D21').Value = '0.0₃0945'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('D22').Value = '2.66'
$ws.Range('E22').Value = '  +13.29%  '
$ws.Range('D23').Value = '70.75'
$ws.Range('E23').Value = '  -0.78%  '
$ws.Range('D24').Value = '249.42'
$ws.Range('E24').Value = '  +1.42%  '
$ws.Range('E25').Value = '  +3.58%  '
$ws.Range('D26').Value = '26.16'
$ws.Range('E26').Value = '  +0.85%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '2.30'
$ws.Range('E28').Value = '  +10.94%  '
$ws.Range('E29').Value = '  +3.59%  '
$ws.Range('E30').Value = '  +5.67%  '
$ws.Range('E31').Value = '  +7.54%  '
$ws.Range('D32').Value = '50.36'
$ws.Range('E32').Value = '  +2.40%  '
$ws.Range('D33').Value = '19.98'
$ws.Range('E33').Value = '  +0.55%  '
$ws.Range('E34').Value = '  +2.01%  '
$ws.Range('E35').Value = '  +4.91%  '
$ws.Range('E36').Value = '  +0.24%  '
$ws.Range('E37').Value = '  +5.96%  '
$ws.Range('D38').Value = '4.76'
$ws.Range('E38').Value = '  +5.01%  '
$ws.Range('E39').Value = '  +3.18%  '
$ws.Range('E40').Value = '  +1.75%  '
$ws.Range('D41').Value = '122.28'
$ws.Range('E41').Value = '  -3.25%  '
$ws.Range('E42').Value = '  -0.93%  '
$ws.Range('D43').Value = '21.34'
$ws.Range('E43').Value = '  +1.38%  '
$ws.Range('D45').Value = '2.006.45'
$ws.Range('E45').Value = '  +2.59%  '
$ws.Range('E46').Value = '  +5.26%  '
$ws.Range('E48').Value = '  -3.34%  '
$ws.Range('D49').Value = '9.08'
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('E50').Value = '  +6.64%  '
$ws.Range('D51').Value = '78.50'
$ws.Range('E51').Value = '  +1.30%  '

# Restore the original (default/general) formatting so no stray cell styles
# are introduced.
$ws.Range("D2:E51").ClearFormats()
